$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q4").Font.Bold = $ws.Range("P4").Font.Bold
$ws.Range("Q4").Font.Name = $ws.Range("P4").Font.Name
$ws.Range("Q4").Font.Size = $ws.Range("P4").Font.Size
$ws.Range("Q4").HorizontalAlignment = $ws.Range("P4").HorizontalAlignment
$ws.Range("Q4").VerticalAlignment = $ws.Range("P4").VerticalAlignment
$ws.Range("Q4").WrapText = $ws.Range("P4").WrapText
$ws.Range("Q4").Borders.LineStyle = $ws.Range("P4").Borders.LineStyle
$ws.Range("Q4").Value = 2020
